# Auto-generated: apply cryptos price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.843.19"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.647.74"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +1.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.871.99"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "1.645.93"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "26.833.91"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("E22").Value = "  +11.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "1.283.09"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.538"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "1.797.21"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  -0.76%  "
